# TC04 & TC13 excel data changes
#
# TC04 (HRM_Attendance): the "In Time" (col F) and "Out Time" (col I) values
# were entered swapped for the 5 sample rows -- fix by swapping them back.
#
# TC13 (Demo_PlaceOrder): rename the generic header labels in row 1
# (A1:D1) to the test-case-specific labels used on this sheet.

$wb = $excel.ActiveWorkbook

# --- TC04: HRM_Attendance - swap "In Time" / "Out Time" values ------------
$wsAttendance = $wb.Worksheets.Item("HRM_Attendance")

for ($r = 2; $r -le 6; $r++) {
    $inTimeCell  = $wsAttendance.Cells.Item($r, 6)   # column F - " In Time"
    $outTimeCell = $wsAttendance.Cells.Item($r, 9)   # column I - " Out Time"

    $inTimeText  = $inTimeCell.Text
    $outTimeText = $outTimeCell.Text

    # Leading "'" keeps these as text (quote-prefixed), matching how the
    # original values were stored, instead of being reinterpreted as times.
    $inTimeCell.Value  = "'" + $outTimeText
    $outTimeCell.Value = "'" + $inTimeText
}

$wsAttendance.Activate()
$wsAttendance.Range("E2").Select()

# --- TC13: Demo_PlaceOrder - rename header row -----------------------------
$wsPlaceOrder = $wb.Worksheets.Item("Demo_PlaceOrder")

$wsPlaceOrder.Range("A1").Value = "Test Case_ID"
$wsPlaceOrder.Range("B1").Value = "Test Case_Name"
$wsPlaceOrder.Range("C1").Value = "User_Name"
$wsPlaceOrder.Range("D1").Value = "Password_DWS"

$wsPlaceOrder.Range("C2").Select()
